# Insert a new row for an additional faculty member
# ("1176388 - Luiz Tadeu Fernandes Eleno") right after the existing
# "519033 - Carlos Yujiro Shigue" row (row 13), under the "Docentes
# responsáveis:" section. Inserting the row shifts every subsequent
# row (the syllabus/requirements rows) down by one, matching the target
# diff (dimension grows from A1:C44 to A1:C45).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 14; this pushes old rows 14..44
# down to 15..45 and copies formatting from the row above (row 13),
# which already carries the B/C column styles we need.
$ws.Rows("14").Insert()

# Populate the new row's faculty-name cells (columns B and C), mirroring
# how the existing "519033 - Carlos Yujiro Shigue" row (B13/C13) is laid
# out -- same text duplicated in both columns.
$ws.Range("B14").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C14").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
